$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new price records for "Vega Monumental Concepción"
# (2022-10-21 / serial 44855) for Pimiento "Zafiro rojo" and "Zafiro verde".
# They are inserted right before the existing row 299, pushing the prior
# rows 299-389 down to 301-391 (dimension grows from R389 to R391).
$ws.Rows("299:300").Insert()

# Row 299: new record - Zafiro rojo
$ws.Range("A299").Value = 11
$ws.Range("B299").Value = "Vega Monumental Concepción"
$ws.Range("C299").Value = "Bíobío"
$ws.Range("D299").Value = 44855
$ws.Range("E299").Value = 8
$ws.Range("F299").Value = 100112002
$ws.Range("G299").Value = "Pimiento"
$ws.Range("H299").Value = "Zafiro rojo"
$ws.Range("I299").Value = "Primera"
$ws.Range("J299").Value = 100
$ws.Range("K299").Value = 12000
$ws.Range("L299").Value = 13000
$ws.Range("M299").Value = 12500
$ws.Range("N299").Value = "`$/caja 15 kilos"
$ws.Range("O299").Value = "Región de Arica y Parinacota"
$ws.Range("P299").Value = 833
$ws.Range("Q299").Value = 15
$ws.Range("R299").Value = "Hortaliza"

# Row 300: new record - Zafiro verde
$ws.Range("A300").Value = 11
$ws.Range("B300").Value = "Vega Monumental Concepción"
$ws.Range("C300").Value = "Bíobío"
$ws.Range("D300").Value = 44855
$ws.Range("E300").Value = 8
$ws.Range("F300").Value = 100112002
$ws.Range("G300").Value = "Pimiento"
$ws.Range("H300").Value = "Zafiro verde"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 100
$ws.Range("K300").Value = 11000
$ws.Range("L300").Value = 12000
$ws.Range("M300").Value = 11500
$ws.Range("N300").Value = "`$/caja 15 kilos"
$ws.Range("O300").Value = "Región de Arica y Parinacota"
$ws.Range("P300").Value = 767
$ws.Range("Q300").Value = 15
$ws.Range("R300").Value = "Hortaliza"
